$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (18) down to the two new rows
# so the new cells pick up the same style indices (date format, fills, etc.)
$ws.Range("A18").Copy()
$ws.Range("A19:A20").PasteSpecial(-4122)

$ws.Range("B18").Copy()
$ws.Range("B19:B20").PasteSpecial(-4122)

$ws.Range("C18").Copy()
$ws.Range("C19:C20").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Row 19: Count of an element in array
$ws.Range("A19").Value = (Get-Date -Year 2025 -Month 11 -Day 7).Date
$ws.Range("B19").Value = "Count of an element in array"
$ws.Range("C19").Value = "BINARY SEARCH"
$ws.Range("D19").Value = "ADITYA VERMA"

# Row 20: Number of times sorted array is rotated
$ws.Range("A20").Value = (Get-Date -Year 2025 -Month 11 -Day 7).Date
$ws.Range("B20").Value = "Number of times sorted array is rotated"
$ws.Range("C20").Value = "BINARY SEARCH"
$ws.Range("D20").Value = "ADITYA VERMA"

$ws.Range("A21").Select()
